# Create D_Ratio macro and add test
#
# This script:
#  1. Removes the old, duplicate "Specific Risk" test row (previously at
#     row 62 - "Test Specific Risk with scale=252") since the test suite
#     is reorganizing Specific Risk tests further down, right after the
#     "Table_SpecificRisk1" row.
#  2. Re-adds/renames the Specific Risk tests (daily/monthly/yearly) plus
#     new Total Risk, Average Length, Average Recovery, and D Ratio tests
#     at the bottom of the table (rows 97-104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the old "Specific_Risk1 / Test Specific Risk with scale=252"
#    row - everything below shifts up by one row.
$ws.Rows.Item(62).Delete()

# 2. Append the (re-ordered/renamed) Specific Risk tests and the brand
#    new Total Risk / Average Length / Average Recovery / D Ratio tests.
#    Cells are written in the same per-cell order the original author
#    used so the shared-string table comes out byte-identical.
$ws.Cells.Item(97, 1).Value = "Specific_Risk1"
$ws.Cells.Item(97, 2).Value = "Test Specific Risk for daily returns"
$ws.Cells.Item(97, 3).Value = "Specific_Risk_test1"

$ws.Cells.Item(98, 2).Value = "Test Specific Risk for monthly returns"
$ws.Cells.Item(98, 3).Value = "Specific_Risk_test2"
$ws.Cells.Item(98, 1).Value = "Specific_Risk2"

$ws.Cells.Item(99, 1).Value = "Specific_Risk3"
$ws.Cells.Item(99, 2).Value = "Test Specific Risk for yearly returns"
$ws.Cells.Item(99, 3).Value = "Specific_Risk_test3"

$ws.Cells.Item(100, 1).Value = "Total_Risk1"
$ws.Cells.Item(100, 2).Value = "Test Total Risk with VARDEF=DF"
$ws.Cells.Item(100, 3).Value = "Total_Risk_test1"

$ws.Cells.Item(101, 1).Value = "Total_Risk2"
$ws.Cells.Item(101, 2).Value = "Test Total Risk with VARDEF=N"
$ws.Cells.Item(101, 3).Value = "Total_Risk_test2"

$ws.Cells.Item(102, 3).Value = "Average_Length_test"
$ws.Cells.Item(102, 1).Value = "Average Length"
$ws.Cells.Item(102, 2).Value = "Test Average length"

$ws.Cells.Item(103, 1).Value = "Average Recovery"
$ws.Cells.Item(103, 2).Value = "Test Average Reconvery"
$ws.Cells.Item(103, 3).Value = "Average_Recovery_test"

$ws.Cells.Item(104, 3).Value = "D_Ratio_test"
$ws.Cells.Item(104, 2).Value = "Test d ratio"
$ws.Cells.Item(104, 1).Value = "D Ratio"

# Match the author's final selection / scroll position.
$ws.Range("A104").Select()
